# "added pstools, update evidence"
# Inserts a new "PsExec" slide (title + content placeholder) right before
# the existing "References" slide, so the deck order becomes:
#   1 BlackCat commands
#   2 Detailed commands
#   3 T1490 Inhibit System Recovery
#   4 (image slide)
#   5 PsExec               <-- new
#   6 References
#
# ppLayoutText (=2) is the classic "Title and Content" autolayout, matching
# the slideLayout2.xml ("Title and Content") used by the References slide.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(5, 2)

# --- Title placeholder -----------------------------------------------
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "PsExec"
$title.LanguageID = "en-AU"

# --- Body / content placeholder ---------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "BlackCat"
$body.LanguageID = "en-AU"

$run = $body.InsertAfter(" will have ")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter("PsExec")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter(" bundled into the executable however ")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter("PsExec")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter(" can be download from ")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter("here")
$run.LanguageID = "en-AU"

$run = $run.InsertAfter(" .")
$run.LanguageID = "en-AU"

# Scope the hyperlink to just the "here" word. Doing this via a fresh
# Characters() lookup (after the full string exists) keeps the hyperlink
# from bleeding onto the sibling runs (InsertAfter().ActionSettings would
# otherwise stamp the hlinkClick across the whole paragraph).
$fullText = $body.Text
$hereStart = $fullText.IndexOf("here") + 1
$hereLink = $body.Characters($hereStart, 4)
$hereLink.ActionSettings.Item(1).Hyperlink.Address = "https://learn.microsoft.com/en-us/sysinternals/downloads/psexec"
